# Selenium_Automation ECommerce TestData.xlsx update:
# Adds three new page-test sheets (OrderShippingPageTest, OrderPaymentPageTest,
# OrderConfirmationPageTest) after the existing OrderDetailsPageTest sheet,
# and updates the active-sheet/selection bookkeeping to match.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) OrderShippingPageTest - a duplicate of OrderDetailsPageTest (sheet5),
#    used as the starting point for the new shipping-page test data.
# ------------------------------------------------------------------
$orderDetails = $wb.Worksheets.Item("OrderDetailsPageTest")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$orderDetails.Copy($null, $lastSheet) | Out-Null

$shipping = $wb.Worksheets.Item($wb.Worksheets.Count)
$shipping.Name = "OrderShippingPageTest"

# ------------------------------------------------------------------
# 2) OrderPaymentPageTest - brand-new sheet with payment summary data.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$payment = $wb.Worksheets.Add($null, $lastSheet)
$payment.Name = "OrderPaymentPageTest"

$payment.Range("A1").Value = "PaymentMode"
$payment.Range("B1").Value = "PaymentText"
$payment.Range("C1").Value = "TotalAmountOfProduct"

$payment.Range("A2:C2").NumberFormat = "@"
$payment.Range("A2").Value = "CHECK PAYMENT"
$payment.Range("B2").Value = "You have chosen to pay by check. Here is a short summary of your order:"
$payment.Range("C2").Value = "`$18.51"

$payment.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# 3) OrderConfirmationPageTest - brand-new sheet with confirmation data.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$confirmation = $wb.Worksheets.Add($null, $lastSheet)
$confirmation.Name = "OrderConfirmationPageTest"

$confirmation.Range("A1:D2").NumberFormat = "@"

$confirmation.Range("A1").Value = "OrderConfimationPageLabel"
$confirmation.Range("B1").Value = "OrderConfimationSuccessMessage"
$confirmation.Range("C1").Value = "OrderConfirmationText"
$confirmation.Range("D1").Value = "TotalAmountOfProduct"

$confirmation.Range("A2").Value = "ORDER CONFIRMATION"
$confirmation.Range("B2").Value = "Your order on My Store is complete."
$confirmation.Range("C2").Value = "Your order will be sent as soon as we receive your payment."
$confirmation.Range("D2").Value = "`$18.51"

$confirmation.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# 4) Selections / active-sheet bookkeeping to match the saved workbook
#    view state: OrderDetailsPageTest and OrderPaymentPageTest end up
#    fully selected (Ctrl+A), OrderShippingPageTest has C14 selected,
#    OrderConfirmationPageTest ends up the active tab with B11 selected.
# ------------------------------------------------------------------
$orderDetails.Activate() | Out-Null
$orderDetails.Cells.Select() | Out-Null

$shipping.Activate() | Out-Null
$shipping.Range("C14").Select() | Out-Null

$payment.Activate() | Out-Null
$payment.Cells.Select() | Out-Null

$confirmation.Activate() | Out-Null
$confirmation.Range("B11").Select() | Out-Null
